$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price / Volume(1h) refresh for most rows (simple cell-value edits).
$ws.Range("D2").Value = "30.381.73"
$ws.Range("E2").Value = "  +2.23%  "
$ws.Range("D3").Value = "2.093.66"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  -0.79%  "
$ws.Range("D5").Value = "343.00"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").Value = "0.5231"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("D8").Value = "0.4420"
$ws.Range("E8").Value = "  +0.98%  "
$ws.Range("D9").Value = "54.57"
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("D10").Value = "0.09313"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "24.79"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "8.598"
$ws.Range("E13").Value = "  +3.25%  "
$ws.Range("D14").Value = "6.907"
$ws.Range("E14").Value = "  +2.47%  "
$ws.Range("D15").Value = "2.089.49"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "101.66"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").Value = "21.16"
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("D20").Value = "0.06674"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "6.332"
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").Value = "30.395.03"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("D24").Value = "12.52"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "2.317"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").Value = "163.12"
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("D28").Value = "2.505"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").Value = "133.15"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "1.135"
$ws.Range("E30").Value = "  +0.50%  "

# Rows 31/32 swapped rank order (ARBITRUM <-> Stellar) with refreshed data.
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.1047"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").Value = "1.656"
$ws.Range("E32").Value = "  +0.49%  "

# Remaining rows.
$ws.Range("D33").Value = "6.826"
$ws.Range("E33").Value = "  +9.64%  "
$ws.Range("D34").Value = "6.255"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("D35").Value = "3.856"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("D36").Value = "10.12"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "0.02631"
$ws.Range("E37").Value = "  +2.37%  "
$ws.Range("D38").Value = "0.06839"
$ws.Range("E38").Value = "  +2.30%  "
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("E40").Value = "  +1.18%  "
$ws.Range("D41").Value = "1.340"
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").Value = "0.2212"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("E43").Value = "  +1.90%  "
$ws.Range("D44").Value = "14.42"
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").Value = "1.382"
$ws.Range("E47").Value = "  +18.98%  "
$ws.Range("D48").Value = "3.634"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").Value = "1.237"
$ws.Range("E49").Value = "  +10.84%  "
$ws.Range("D50").Value = "0.00000000348"
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("D51").Value = "1.216"
$ws.Range("E51").Value = "  +0.00%  "
